$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Formula = '="27.062.02"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E2")
$c.Formula = '="  +0.43%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D3")
$c.Formula = '="1.562.54"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E3")
$c.Formula = '="  +0.51%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E4")
$c.Formula = '="  +0.28%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D5")
$c.Formula = '="210.48"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E5")
$c.Formula = '="  +1.66%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E6")
$c.Formula = '="  +0.26%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E7")
$c.Formula = '="  +0.39%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D8")
$c.Formula = '="21.89"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E8")
$c.Formula = '="  -0.70%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E9")
$c.Formula = '="  -0.11%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E10")
$c.Formula = '="  +0.31%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D11")
$c.Formula = '="0.0860"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E11")
$c.Formula = '="  +0.34%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D12")
$c.Formula = '="1.784.23"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E12")
$c.Formula = '="  +0.49%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D13")
$c.Formula = '="1.553.72"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E13")
$c.Formula = '="  -0.10%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D14")
$c.Formula = '="3.76"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E14")
$c.Formula = '="  +0.18%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E15")
$c.Formula = '="  -0.45%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D16")
$c.Formula = '="27.061.03"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E16")
$c.Formula = '="  +0.43%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D17")
$c.Formula = '="61.92"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E17")
$c.Formula = '="  +0.37%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D18")
$c.Formula = '="0.0₃0700"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E18")
$c.Formula = '="  -0.94%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D19")
$c.Formula = '="214.84"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E19")
$c.Formula = '="  -1.25%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D20")
$c.Formula = '="7.35"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E20")
$c.Formula = '="  +0.55%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E21")
$c.Formula = '="  +0.37%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E22")
$c.Formula = '="  +0.90%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E23")
$c.Formula = '="  -0.11%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E24")
$c.Formula = '="  -0.21%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D25")
$c.Formula = '="153.79"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E25")
$c.Formula = '="  +0.21%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E26")
$c.Formula = '="  -0.67%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D27")
$c.Formula = '="15.01"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E27")
$c.Formula = '="  +0.07%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E28")
$c.Formula = '="  +1.31%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E29")
$c.Formula = '="  +0.40%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E30")
$c.Formula = '="  +4.15%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D31")
$c.Formula = '="0.0471"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E31")
$c.Formula = '="  +0.24%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E32")
$c.Formula = '="  +0.29%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E33")
$c.Formula = '="  +2.00%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D34")
$c.Formula = '="1.430.66"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E34")
$c.Formula = '="  +0.90%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E35")
$c.Formula = '="  +0.16%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E36")
$c.Formula = '="  -0.53%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E37")
$c.Formula = '="  +1.81%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E38")
$c.Formula = '="  +0.80%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D39")
$c.Formula = '="0.530"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E39")
$c.Formula = '="  +0.49%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E40")
$c.Formula = '="  +2.75%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E41")
$c.Formula = '="  -0.07%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E42")
$c.Formula = '="  +0.43%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D43")
$c.Formula = '="2.34"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E43")
$c.Formula = '="  +1.24%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E44")
$c.Formula = '="  +0.13%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D45")
$c.Formula = '="64.31"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E45")
$c.Formula = '="  -0.36%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D47")
$c.Formula = '="1.703.10"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D48")
$c.Formula = '="85.90"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E48")
$c.Formula = '="  -1.52%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E49")
$c.Formula = '="  +2.40%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D50")
$c.Formula = '="0.0517"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E50")
$c.Formula = '="  -0.70%  "'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("D51")
$c.Formula = '="0.0957"'
$c.Copy()
$c.PasteSpecial(-4163)
$c = $ws.Range("E51")
$c.Formula = '="  -0.34%  "'
$c.Copy()
$c.PasteSpecial(-4163)
